$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nov")

# Merge "Test" (C4) and "Complete" (D4) into a single value "TestComplete" in C4,
# then clear out D4 since that column's value is no longer present.
$ws.Range("C4").Value = "TestComplete"
$ws.Range("D4").ClearContents()

# Update the active selection on the sheet to C5.
$ws.Activate()
$ws.Range("C5").Select()
